# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E27) currently lists the contribution
# periods in descending order (2205 .. 2106). This update refreshes the
# database so the periods run in ascending chronological order
# (2106 .. 2205), i.e. part 1 of the new "Estado de Cuenta" periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2106","2107","2108","2109","2110","2111","2112","2201","2202","2203","2204","2205")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}
